$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect so the cells below can be written.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer cell (A10)
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-14 for illustrative purposes only and are subject to change."

# Update weight / percent-change figures for rows 2-7
$ws.Range("D2").Value = 0.4977034582930538
$ws.Range("E2").Value = 0.001152959262106235

$ws.Range("D3").Value = 0.3248307044347181
$ws.Range("E3").Value = -0.0001906759462294572

$ws.Range("D4").Value = 0.09106923387077187
$ws.Range("E4").Value = 0.008168028004667427

$ws.Range("D5").Value = 0.058586102704357
$ws.Range("E5").Value = 0.001596169193934571

$ws.Range("D6").Value = 0.02781050069709939
$ws.Range("E6").Value = -0.01532784558614808

$ws.Range("E7").Value = 0.0009229887346795529

# Restore sheet protection (matching the original protected state).
$ws.Protect("D382", $true, $true, $true)
